# Daily attendance processing - swap the "Recorded By" ordering in column G
# from "System, <email>" to "<email>, System" for every affected row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace(
    "System, dnasr281@gmail.com",
    "dnasr281@gmail.com, System",
    [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole
)
